# study.xlsx update — "Add files via upload"
#
# Adds five new daily log rows (9-13) to the study tracker: the date
# number in column A, the "time" note in the merged B:F cell, and the
# "content" note in the merged G:N cell (left blank for row 10/G6, which
# the author didn't fill in). Also moves the active selection to the
# newly-entered row (B14:F14) and resizes the workbook window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the saved window chrome to match the re-uploaded file.
$win = $excel.ActiveWindow
$win.Width  = 12800
$win.Height = 11630

# Row 5 (day 9)
$ws.Range("A5").Value = 9
$ws.Range("B5").Value = "正常晚自习时间，回来后也学了一会儿"
$ws.Range("G5").Value = "高数，c语言"

# Row 6 (day 10) — only the time column was filled in ("same as above")
$ws.Range("A6").Value = 10
$ws.Range("B6").Value = "同上"

# Row 7 (day 11)
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "正常晚自习时间"
$ws.Range("G7").Value = "（下周月考，这几天事情也有些多，忘了传，不过进度没落下）"

# Row 8 (day 12)
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "4-5：30，正常晚自习"
$ws.Range("G8").Value = "线代，高数，C语言"

# Row 9 (day 13)
$ws.Range("A9").Value = 13
$ws.Range("B9").Value = "4-5：30，6：30-9"
$ws.Range("G9").Value = "高数，C语言"

# Leave the selection on the next empty row, as in the source file.
$null = $ws.Range("B14:F14").Select()
